$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Best-effort: restore window position / scroll (may not be serialized by this engine) ---
$win = $excel.ActiveWindow
$win.Left = 6690
$win.Top = 3015
$win.ScrollRow = 21
$win.ScrollColumn = 1

# --- Resize column B (User Story) from 111.71 chars to 59 chars stored width ---
$ws.Columns.Item(2).ColumnWidth = 58.1666666666667

# --- Seed the 5 new lookup values into the shared string table in the exact order
#     they first appear in the target workbook (To Do, High, Urgent, Normal, Low),
#     using a scratch range far outside the used area, then clear the scratch cells. ---
$ws.Range("Z100").Value = "To Do"
$ws.Range("Z101").Value = "High"
$ws.Range("Z102").Value = "Urgent"
$ws.Range("Z103").Value = "Normal"
$ws.Range("Z104").Value = "Low"

# --- Fill in Sprint Number (C), Story Priority (D), Story Status (E), Story Points (F) ---
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "Urgent"
$ws.Range("E2").Value = "To Do"
$ws.Range("F2").Value = 3

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Urgent"
$ws.Range("E3").Value = "In Progress"
$ws.Range("F3").Value = 3

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "Urgent"
$ws.Range("E4").Value = "In Progress"
$ws.Range("F4").Value = 3

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "Normal"
$ws.Range("E5").Value = "To Do"
$ws.Range("F5").Value = 10

$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Normal"
$ws.Range("E6").Value = "To Do"
$ws.Range("F6").Value = 10

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "High"
$ws.Range("E7").Value = "To Do"
$ws.Range("F7").Value = 3

$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "Low"
$ws.Range("E8").Value = "To Do"
$ws.Range("F8").Value = 5

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "Urgent"
$ws.Range("E9").Value = "In Progress"
$ws.Range("F9").Value = 5

$ws.Range("D10").Value = "Normal"
$ws.Range("E10").Value = "To Do"
$ws.Range("F10").Value = 20

$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "High"
$ws.Range("E11").Value = "To Do"
$ws.Range("F11").Value = 5

$ws.Range("C12").Value = 3
$ws.Range("D12").Value = "High"
$ws.Range("E12").Value = "To Do"
$ws.Range("F12").Value = 3

$ws.Range("C13").Value = 4
$ws.Range("D13").Value = "Normal"
$ws.Range("E13").Value = "To Do"
$ws.Range("F13").Value = 5

$ws.Range("C14").Value = 4
$ws.Range("D14").Value = "Normal"
$ws.Range("E14").Value = "To Do"
$ws.Range("F14").Value = 8

$ws.Range("C15").Value = 4
$ws.Range("D15").Value = "Normal"
$ws.Range("E15").Value = "To Do"
$ws.Range("F15").Value = 20

$ws.Range("C16").Value = 4
$ws.Range("D16").Value = "High"
$ws.Range("E16").Value = "To Do"
$ws.Range("F16").Value = 10

$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "High"
$ws.Range("E17").Value = "To Do"
$ws.Range("F17").Value = 5

$ws.Range("C18").Value = 3
$ws.Range("D18").Value = "Urgent"
$ws.Range("E18").Value = "To Do"
$ws.Range("F18").Value = 8

$ws.Range("C19").Value = 5
$ws.Range("D19").Value = "Normal"
$ws.Range("E19").Value = "To Do"
$ws.Range("F19").Value = 20

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "Urgent"
$ws.Range("E20").Value = "To Do"
$ws.Range("F20").Value = 8

$ws.Range("C21").Value = 2
$ws.Range("D21").Value = "Urgent"
$ws.Range("E21").Value = "To Do"
$ws.Range("F21").Value = 8

$ws.Range("C22").Value = 4
$ws.Range("D22").Value = "Normal"
$ws.Range("E22").Value = "To Do"
$ws.Range("F22").Value = 20

$ws.Range("C23").Value = 3
$ws.Range("D23").Value = "High"
$ws.Range("E23").Value = "To Do"
$ws.Range("F23").Value = 20

$ws.Range("C24").Value = 2
$ws.Range("D24").Value = "Urgent"
$ws.Range("E24").Value = "To Do"
$ws.Range("F24").Value = 5

# --- Remove the scratch seed cells now that the real cells reference the strings ---
$ws.Range("Z100:Z104").ClearContents() | Out-Null

# --- Row heights increased where the narrower column B causes extra text wrapping ---
$ws.Rows.Item(7).RowHeight = 47.25
$ws.Rows.Item(16).RowHeight = 47.25
$ws.Rows.Item(18).RowHeight = 47.25
$ws.Rows.Item(19).RowHeight = 63
$ws.Rows.Item(20).RowHeight = 47.25
$ws.Rows.Item(21).RowHeight = 63
$ws.Rows.Item(22).RowHeight = 47.25
$ws.Rows.Item(23).RowHeight = 63
$ws.Rows.Item(24).RowHeight = 47.25

# --- Update selection / scroll position ---
$ws.Range("C6").Select() | Out-Null
